# Add the "stim details" block at the bottom of the sheet, and fill in the
# pair_kind ("generic") values for the practice pairs in column J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice pairs (rows 2-5) didn't have a pair_kind (column J) value yet.
$ws.Range("J2:J5").Value = "generic"

# New "stim details" mini-table starting a few rows below the existing data.
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$stimDetails = @(
    @(6, "video"),
    @(6, "video"),
    @(7, "video"),
    @(7, "video"),
    @(6, "audio"),
    @(6, "audio"),
    @(7, "audio"),
    @(7, "audio")
)

$row = 29
foreach ($entry in $stimDetails) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
